$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.923.88"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "3.523.77"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'597.01"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "'134.05"
$ws.Range("E6").Value = "  -2.00%  "
$ws.Range("D7").Value = "3.520.92"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.496"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").Value = "'0.124"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").Value = "'7.14"
$ws.Range("E11").Value = "  +3.52%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "4.128.43"
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("D14").Value = "'27.27"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "3.521.45"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "64.953.24"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'9.98"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").Value = "'14.43"
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("D21").Value = "'5.69"
$ws.Range("E21").Value = "  -2.35%  "
$ws.Range("D22").Value = "'390.76"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").Value = "3.666.46"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "'74.26"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").Value = "'1.61"
$ws.Range("E28").Value = "  +19.21%  "
$ws.Range("D29").Value = "'7.75"
$ws.Range("E29").Value = "  +1.17%  "
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("E31").Value = "  +1.59%  "
$ws.Range("D32").Value = "'8.38"
$ws.Range("E32").Value = "  +2.55%  "
$ws.Range("D33").Value = "3.527.94"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("D34").Value = "'24.12"
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("E37").Value = "  +5.89%  "
$ws.Range("E38").Value = "  +2.03%  "
$ws.Range("D39").Value = "'168.65"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").Value = "'6.83"
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("D41").Value = "'0.0820"
$ws.Range("E41").Value = "  +3.01%  "
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "'42.67"
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").Value = "'1.24"
$ws.Range("E44").Value = "  +3.69%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").Value = "'25.28"
$ws.Range("E46").Value = "  -4.48%  "
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("D49").Value = "'6.91"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("D50").Value = "2.401.47"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").Value = "'0.897"
$ws.Range("E51").Value = "  +6.60%  "
